# Apply "repull data, push all data, mean calculation" updates to column F (dSF)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> new value for column F
$updates = @{
    2  = -7
    5  = -8
    6  = -7
    8  = 7
    9  = -2
    11 = -2
    13 = -7
    14 = -4
    15 = 4
    16 = 10
    17 = -4
    19 = 9
    24 = -5
    25 = 1
    28 = -10
    32 = 1
    33 = 2
    34 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
